$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-10 with corrected simulation values ---
# Row 2
$ws.Cells.Item(2, 4).Value2 = 0.9755265073994462
$ws.Cells.Item(2, 5).Value2 = 0.7894726959413937
$ws.Cells.Item(2, 6).Value2 = 0.934335133894456
$ws.Cells.Item(2, 7).Value2 = 0.8125346136630779
$ws.Cells.Item(2, 8).Value2 = 0.7802116002846456
$ws.Cells.Item(2, 9).Value2 = 0.7893851990814973
$ws.Cells.Item(2, 10).Value2 = 0.8098178552922931
$ws.Cells.Item(2, 11).Value2 = 0.881585251661231
$ws.Cells.Item(2, 12).Value2 = 0.9960483123262905
$ws.Cells.Item(2, 13).Value2 = 0.6453029786796329
$ws.Cells.Item(2, 14).Value2 = 0.791120953156803
$ws.Cells.Item(2, 15).Value2 = 0.9475850065585952
$ws.Cells.Item(2, 16).Value2 = 0.9071149855068927
$ws.Cells.Item(2, 17).Value2 = 0.9331797801905444
$ws.Cells.Item(2, 18).Value2 = 0.781219349366597
$ws.Cells.Item(2, 19).Value2 = 0.9118501324214733
$ws.Cells.Item(2, 20).Value2 = 0.9336362698648814
$ws.Cells.Item(2, 21).Value2 = 0.8071756365407339
$ws.Cells.Item(2, 22).Value2 = 0.93080314965507
$ws.Cells.Item(2, 23).Value2 = 0.9057297649285171
$ws.Cells.Item(2, 24).Value2 = 0.8606375026371246

# Row 3
$ws.Cells.Item(3, 4).Value2 = 0.7399255225012077
$ws.Cells.Item(3, 5).Value2 = 0.650079016375349
$ws.Cells.Item(3, 6).Value2 = 0.8959836739881956
$ws.Cells.Item(3, 7).Value2 = 0.4165861382681959
$ws.Cells.Item(3, 8).Value2 = 0.8944323322249178
$ws.Cells.Item(3, 9).Value2 = 0.8421052631578947
$ws.Cells.Item(3, 10).Value2 = 0.5346250261368419
$ws.Cells.Item(3, 11).Value2 = 0.902944466371997
$ws.Cells.Item(3, 12).Value2 = 0.8980942128756545
$ws.Cells.Item(3, 13).Value2 = 0.3004127618280008
$ws.Cells.Item(3, 14).Value2 = 0.7494643545870795
$ws.Cells.Item(3, 15).Value2 = 0.8435231622694981
$ws.Cells.Item(3, 16).Value2 = 0.8385332617800386
$ws.Cells.Item(3, 17).Value2 = 0.820342587241271
$ws.Cells.Item(3, 18).Value2 = 0.3911480060021357
$ws.Cells.Item(3, 19).Value2 = 0.9166260621878772
$ws.Cells.Item(3, 20).Value2 = 0.9126032692761412
$ws.Cells.Item(3, 21).Value2 = 0.681690133074805
$ws.Cells.Item(3, 22).Value2 = 0.9268814852801942
$ws.Cells.Item(3, 23).Value2 = 0.7592705648032437
$ws.Cells.Item(3, 24).Value2 = 0.8287873670856611

# Row 4
$ws.Cells.Item(4, 4).Value2 = 0.8330369197824673
$ws.Cells.Item(4, 5).Value2 = 0.6012296051001795
$ws.Cells.Item(4, 6).Value2 = 0.6877957928437578
$ws.Cells.Item(4, 7).Value2 = 0.927379019010212
$ws.Cells.Item(4, 8).Value2 = 0.8942653144281929
$ws.Cells.Item(4, 9).Value2 = 0.842104439878995
$ws.Cells.Item(4, 10).Value2 = 0.8058722812113154
$ws.Cells.Item(4, 11).Value2 = 0.8730899510995821
$ws.Cells.Item(4, 12).Value2 = 0.9869755930323554
$ws.Cells.Item(4, 13).Value2 = 0.3232019907052007
$ws.Cells.Item(4, 14).Value2 = 0.7406080125687596
$ws.Cells.Item(4, 15).Value2 = 0.916841899110767
$ws.Cells.Item(4, 16).Value2 = 0.862090156474621
$ws.Cells.Item(4, 17).Value2 = 0.9329219691301992
$ws.Cells.Item(4, 18).Value2 = 0.753571629303206
$ws.Cells.Item(4, 19).Value2 = 0.9289836028266324
$ws.Cells.Item(4, 20).Value2 = 0.5538354070944473
$ws.Cells.Item(4, 21).Value2 = 0.8928625563463907
$ws.Cells.Item(4, 22).Value2 = 0.9292810738774018
$ws.Cells.Item(4, 23).Value2 = 0.8261182155237364
$ws.Cells.Item(4, 24).Value2 = 0.8553252767994508

# Row 5
$ws.Cells.Item(5, 2).Value2 = 'hisditonly'
$ws.Cells.Item(5, 3).Value2 = 'raw'
$ws.Cells.Item(5, 4).Value2 = 0.8300200297589607
$ws.Cells.Item(5, 5).Value2 = 0.6193737702786831
$ws.Cells.Item(5, 6).Value2 = 0.551958888325472
$ws.Cells.Item(5, 7).Value2 = 0.1147762086215552
$ws.Cells.Item(5, 8).Value2 = 0.8370080968595582
$ws.Cells.Item(5, 9).Value2 = 0.8421052631578947
$ws.Cells.Item(5, 10).Value2 = 0.8053916207744045
$ws.Cells.Item(5, 11).Value2 = 0.8590856058479215
$ws.Cells.Item(5, 12).Value2 = 0.8860313827645816
$ws.Cells.Item(5, 13).Value2 = 0.4528544604989494
$ws.Cells.Item(5, 14).Value2 = 0.8245538882612222
$ws.Cells.Item(5, 15).Value2 = 0.9282721941529722
$ws.Cells.Item(5, 16).Value2 = 0.8358425952703649
$ws.Cells.Item(5, 17).Value2 = 0.9030978363569152
$ws.Cells.Item(5, 18).Value2 = 0.5880341337078518
$ws.Cells.Item(5, 19).Value2 = 0.7450095353668735
$ws.Cells.Item(5, 20).Value2 = 0.9487154578808976
$ws.Cells.Item(5, 21).Value2 = 0.8178332978335928
$ws.Cells.Item(5, 22).Value2 = 0.8214321533627195
$ws.Cells.Item(5, 23).Value2 = 0.4926699956846632
$ws.Cells.Item(5, 24).Value2 = 0.740524492408409

# Row 6
$ws.Cells.Item(6, 2).Value2 = 'flowbot'
$ws.Cells.Item(6, 4).Value2 = 0.9893460897462982
$ws.Cells.Item(6, 5).Value2 = 0.7966308621776094
$ws.Cells.Item(6, 6).Value2 = 0.9142945002061762
$ws.Cells.Item(6, 7).Value2 = 0.8057370110753217
$ws.Cells.Item(6, 8).Value2 = 0.9473846568975813
$ws.Cells.Item(6, 9).Value2 = 0.7894736842105263
$ws.Cells.Item(6, 10).Value2 = 0.8034210526670454
$ws.Cells.Item(6, 11).Value2 = 0.8944675794283877
$ws.Cells.Item(6, 12).Value2 = 0.9985735109168359
$ws.Cells.Item(6, 13).Value2 = 0.628598317814225
$ws.Cells.Item(6, 14).Value2 = 0.8974426396949045
$ws.Cells.Item(6, 15).Value2 = 0.9456375257344539
$ws.Cells.Item(6, 16).Value2 = 0.9519444943407248
$ws.Cells.Item(6, 17).Value2 = 0.942644375734421
$ws.Cells.Item(6, 18).Value2 = 0.7796784640513733
$ws.Cells.Item(6, 19).Value2 = 0.8937663946701478
$ws.Cells.Item(6, 20).Value2 = 0.5516436545690855
$ws.Cells.Item(6, 21).Value2 = 0.8767227868629828
$ws.Cells.Item(6, 22).Value2 = 0.9371235708837938
$ws.Cells.Item(6, 23).Value2 = 0.789424868043381
$ws.Cells.Item(6, 24).Value2 = 0.9373001559218095

# Row 7
$ws.Cells.Item(7, 2).Value2 = 'dit'
$ws.Cells.Item(7, 4).Value2 = 0.9450892747086108
$ws.Cells.Item(7, 5).Value2 = 0.6437633304901939
$ws.Cells.Item(7, 6).Value2 = 0.8641820747186533
$ws.Cells.Item(7, 7).Value2 = 0.3767935980948522
$ws.Cells.Item(7, 8).Value2 = 0.9308254275563886
$ws.Cells.Item(7, 9).Value2 = 0.9473684210526315
$ws.Cells.Item(7, 10).Value2 = 0.9283110371107012
$ws.Cells.Item(7, 11).Value2 = 0.8862963151761937
$ws.Cells.Item(7, 12).Value2 = 0.8740837256639908
$ws.Cells.Item(7, 13).Value2 = 0.6463792038533641
$ws.Cells.Item(7, 14).Value2 = 0.8779168287963034
$ws.Cells.Item(7, 15).Value2 = 0.9408511901063883
$ws.Cells.Item(7, 16).Value2 = 0.9123320148265177
$ws.Cells.Item(7, 17).Value2 = 0.9053294353456692
$ws.Cells.Item(7, 18).Value2 = 0.7681093670447995
$ws.Cells.Item(7, 19).Value2 = 0.9158032609625261
$ws.Cells.Item(7, 20).Value2 = 0.938493300704862
$ws.Cells.Item(7, 21).Value2 = 0.7353870272060891
$ws.Cells.Item(7, 22).Value2 = 0.9266557607734368
$ws.Cells.Item(7, 23).Value2 = 0.7124979325785584
$ws.Cells.Item(7, 24).Value2 = 0.8740637685703713

# Row 8
$ws.Cells.Item(8, 2).Value2 = 'pndit'
$ws.Cells.Item(8, 4).Value2 = 0.9802486727122209
$ws.Cells.Item(8, 5).Value2 = 0.8242270995175618
$ws.Cells.Item(8, 6).Value2 = 0.9119628496174405
$ws.Cells.Item(8, 7).Value2 = 0.9343160674474144
$ws.Cells.Item(8, 8).Value2 = 0.9193289248817312
$ws.Cells.Item(8, 9).Value2 = 0.9473684210526315
$ws.Cells.Item(8, 10).Value2 = 0.9397052173340176
$ws.Cells.Item(8, 11).Value2 = 0.8798959836893087
$ws.Cells.Item(8, 12).Value2 = 0.9342260228849276
$ws.Cells.Item(8, 13).Value2 = 0.5776039341722788
$ws.Cells.Item(8, 14).Value2 = 0.800291455099784
$ws.Cells.Item(8, 15).Value2 = 0.939994088863076
$ws.Cells.Item(8, 16).Value2 = 0.9147463553166314
$ws.Cells.Item(8, 17).Value2 = 0.9388001922202354
$ws.Cells.Item(8, 18).Value2 = 0.9513644023661608
$ws.Cells.Item(8, 19).Value2 = 0.9233849975909697
$ws.Cells.Item(8, 20).Value2 = 0.5857930464137202
$ws.Cells.Item(8, 21).Value2 = 0.9199567011628076
$ws.Cells.Item(8, 22).Value2 = 0.929672009486028
$ws.Cells.Item(8, 23).Value2 = 0.7158146313798441
$ws.Cells.Item(8, 24).Value2 = 0.933348180141546

# Row 9
$ws.Cells.Item(9, 2).Value2 = 'hisdit'
$ws.Cells.Item(9, 4).Value2 = 0.929016041851456
$ws.Cells.Item(9, 5).Value2 = 0.8820703340345193
$ws.Cells.Item(9, 6).Value2 = 0.9251687746331324
$ws.Cells.Item(9, 7).Value2 = 0.8121462336601266
$ws.Cells.Item(9, 8).Value2 = 0.9421288875584211
$ws.Cells.Item(9, 9).Value2 = 0.9473684210526315
$ws.Cells.Item(9, 10).Value2 = 0.9367681530347541
$ws.Cells.Item(9, 11).Value2 = 0.9371225908100164
$ws.Cells.Item(9, 12).Value2 = 0.9901238738679399
$ws.Cells.Item(9, 13).Value2 = 0.6309772675465674
$ws.Cells.Item(9, 14).Value2 = 0.915573041892964
$ws.Cells.Item(9, 15).Value2 = 0.9527134496139692
$ws.Cells.Item(9, 16).Value2 = 0.9482375598565721
$ws.Cells.Item(9, 17).Value2 = 0.9338145686777226
$ws.Cells.Item(9, 18).Value2 = 0.7081167992174581
$ws.Cells.Item(9, 19).Value2 = 0.8250852709005912
$ws.Cells.Item(9, 20).Value2 = 0.9258778502021664
$ws.Cells.Item(9, 21).Value2 = 0.8804279106762901
$ws.Cells.Item(9, 22).Value2 = 0.9331765188835016
$ws.Cells.Item(9, 23).Value2 = 0.8691913493962947
$ws.Cells.Item(9, 24).Value2 = 0.9288309946537378

# Row 10
$ws.Cells.Item(10, 2).Value2 = 'pnhisdit'
$ws.Cells.Item(10, 3).Value2 = 'sgp'
$ws.Cells.Item(10, 4).Value2 = 0.9869700660714361
$ws.Cells.Item(10, 5).Value2 = 0.7911992320416115
$ws.Cells.Item(10, 6).Value2 = 0.917535624889571
$ws.Cells.Item(10, 7).Value2 = 0.9257017479344009
$ws.Cells.Item(10, 8).Value2 = 0.9215647119852072
$ws.Cells.Item(10, 9).Value2 = 0.8421052631578947
$ws.Cells.Item(10, 10).Value2 = 0.9274003972909906
$ws.Cells.Item(10, 11).Value2 = 0.8845756533023688
$ws.Cells.Item(10, 12).Value2 = 0.994949148925648
$ws.Cells.Item(10, 13).Value2 = 0.6817895529509904
$ws.Cells.Item(10, 14).Value2 = 0.9273606886408876
$ws.Cells.Item(10, 15).Value2 = 0.9436024651632386
$ws.Cells.Item(10, 16).Value2 = 0.9536160702537589
$ws.Cells.Item(10, 17).Value2 = 0.942188549043224
$ws.Cells.Item(10, 18).Value2 = 0.7779017779003596
$ws.Cells.Item(10, 19).Value2 = 0.9194976791979796
$ws.Cells.Item(10, 20).Value2 = 0.941431051925657
$ws.Cells.Item(10, 21).Value2 = 0.9378828548104982
$ws.Cells.Item(10, 22).Value2 = 0.9334081269107584
$ws.Cells.Item(10, 23).Value2 = 0.8369401539443845
$ws.Cells.Item(10, 24).Value2 = 0.938599445867626

# --- Insert two new data rows (11 and 12), copying formatting from row 10 ---
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 11
$ws.Cells.Item(11, 1).Value2 = 9
$ws.Cells.Item(11, 2).Value2 = 'hisditonly'
$ws.Cells.Item(11, 3).Value2 = 'sgp'
$ws.Cells.Item(11, 4).Value2 = 0.9070990206492752
$ws.Cells.Item(11, 5).Value2 = 0.7074507117110422
$ws.Cells.Item(11, 6).Value2 = 0.8834328832312981
$ws.Cells.Item(11, 7).Value2 = 0.1172589228022622
$ws.Cells.Item(11, 8).Value2 = 0.9290425428288088
$ws.Cells.Item(11, 9).Value2 = 1
$ws.Cells.Item(11, 10).Value2 = 0.9102833447546136
$ws.Cells.Item(11, 11).Value2 = 0.9031071339054562
$ws.Cells.Item(11, 12).Value2 = 0.9419485740396532
$ws.Cells.Item(11, 13).Value2 = 0.619020991173625
$ws.Cells.Item(11, 14).Value2 = 0.8595622980831905
$ws.Cells.Item(11, 15).Value2 = 0.9489195120308755
$ws.Cells.Item(11, 16).Value2 = 0.9149784715701569
$ws.Cells.Item(11, 17).Value2 = 0.9007035115008699
$ws.Cells.Item(11, 18).Value2 = 0.719856285171747
$ws.Cells.Item(11, 19).Value2 = 0.9253519530500142
$ws.Cells.Item(11, 20).Value2 = 0.5977186122970675
$ws.Cells.Item(11, 21).Value2 = 0.9432799508008872
$ws.Cells.Item(11, 22).Value2 = 0.8935923577883212
$ws.Cells.Item(11, 23).Value2 = 0.5477229188435465
$ws.Cells.Item(11, 24).Value2 = 0.7211607363465341

# Row 12
$ws.Cells.Item(12, 1).Value2 = 10
$ws.Cells.Item(12, 4).Value2 = 12
$ws.Cells.Item(12, 5).Value2 = 24
$ws.Cells.Item(12, 6).Value2 = 4
$ws.Cells.Item(12, 7).Value2 = 8
$ws.Cells.Item(12, 8).Value2 = 34
$ws.Cells.Item(12, 9).Value2 = 19
$ws.Cells.Item(12, 10).Value2 = 7
$ws.Cells.Item(12, 11).Value2 = 45
$ws.Cells.Item(12, 12).Value2 = 18
$ws.Cells.Item(12, 13).Value2 = 19
$ws.Cells.Item(12, 14).Value2 = 48
$ws.Cells.Item(12, 15).Value2 = 108
$ws.Cells.Item(12, 16).Value2 = 594
$ws.Cells.Item(12, 17).Value2 = 25
$ws.Cells.Item(12, 18).Value2 = 6
$ws.Cells.Item(12, 19).Value2 = 9
$ws.Cells.Item(12, 20).Value2 = 2
$ws.Cells.Item(12, 21).Value2 = 11
$ws.Cells.Item(12, 22).Value2 = 28
$ws.Cells.Item(12, 23).Value2 = 27
$ws.Cells.Item(12, 24).Value2 = 13

Write-Output "edit applied"